$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
}
